$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full updated data table (rows 2-7, columns A-T), replacing the previous
# 3-row ECs/FAPs/sCs -> ECs-only table with a 7-row ECs/FAPs/sCs -> ECs/sCs table.
$data = @(
    @("ECs",  "Gnas", "Vipr1", "ECs", 3, 1, 136.287657,          408.862971,        0.2628768458810872, 0.2628768458810872, 3, 1,                  1.287673666666667, 3.863021, 0.4047170919281118, 0.4047170919281118, 175.494027010599,  1579.446243095391, 0.1063907526002281, 0.1063907526002281),
    @("ECs",  "Gnas", "Vipr1", "sCs", 3, 1, 136.287657,          408.862971,        0.2628768458810872, 0.2628768458810872, 2, 0.6666666666666666, 1.89399,           5.68197,  0.5952829080718882, 0.5952829080718882, 258.12745948143,   2323.14713533287,  0.1564860932808592, 0.1564860932808592),
    @("FAPs", "Gnas", "Vipr1", "ECs", 3, 1, 143.539174,          430.617522,        0.2768638492442244, 0.2768638492442244, 3, 1,                  1.287673666666667, 3.863021, 0.4047170919281118, 0.4047170919281118, 184.8316144948847, 1663.484530453962, 0.1120515319261457, 0.1120515319261456),
    @("FAPs", "Gnas", "Vipr1", "sCs", 3, 1, 143.539174,          430.617522,        0.2768638492442244, 0.2768638492442244, 2, 0.6666666666666666, 1.89399,           5.68197,  0.5952829080718882, 0.5952829080718882, 271.86176016426,   2446.75584147834,  0.1648123173180788, 0.1648123173180787),
    @("sCs",  "Gnas", "Vipr1", "ECs", 3, 1, 238.6199593333333,   715.859878,        0.4602593048746885, 0.4602593048746884, 3, 1,                  1.287673666666667, 3.863021, 0.4047170919281118, 0.4047170919281118, 307.2646379746042, 2765.381741771438, 0.1862748074017382, 0.1862748074017381),
    @("sCs",  "Gnas", "Vipr1", "sCs", 3, 1, 238.6199593333333,   715.859878,        0.4602593048746885, 0.4602593048746884, 2, 0.6666666666666666, 1.89399,           5.68197,  0.5952829080718882, 0.5952829080718882, 451.9438167777399, 4067.49435099966,  0.2739844974729503, 0.2739844974729503)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $colNum = $c + 1
        $ws.Cells.Item($rowNum, $colNum).Value = $rowVals[$c]
    }
}
